$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H33").Value = 292.23077
$ws.Range("I33").Value = 145.22223
$ws.Range("K33").Value = 145.22223
$ws.Range("M33").Value = 83.77777
$ws.Range("H38").Value = 42.857143
$ws.Range("I38").Value = 33.333332
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 99.999996
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = 272.000004
$ws.Range("N38").Value = -1044
$ws.Range("H64").Value = 5500
$ws.Range("H67").Value = 5500
$ws.Range("H70").Value = 3707.6924
$ws.Range("I70").Value = 3133.6667
$ws.Range("J70").Value = 3879.9
$ws.Range("K70").Value = 9401.000100000001
$ws.Range("L70").Value = 11639.7
$ws.Range("M70").Value = -9131.000100000001
$ws.Range("N70").Value = -12179.7
$ws.Range("H73").Value = 3707.6924
$ws.Range("I73").Value = 3133.6667
$ws.Range("J73").Value = 3879.9
$ws.Range("K73").Value = 9401.000100000001
$ws.Range("L73").Value = 11639.7
$ws.Range("M73").Value = -8465.000100000001
$ws.Range("N73").Value = -13511.7
$ws.Range("H76").Value = 6676.75
$ws.Range("I76").Value = 5651.5
$ws.Range("K76").Value = 5651.5
$ws.Range("M76").Value = -5336.5
$ws.Range("H79").Value = 6676.75
$ws.Range("I79").Value = 5651.5
$ws.Range("K79").Value = 5651.5
$ws.Range("M79").Value = -4559.5
$ws.Range("H112").Value = 2984.3125
$ws.Range("J112").Value = 3541.5833
$ws.Range("L112").Value = 10624.7499
$ws.Range("N112").Value = -12840.7499
$ws.Range("H129").Value = 1924.0834
$ws.Range("I129").Value = 861.7143
$ws.Range("K129").Value = 2585.1429
$ws.Range("M129").Value = 2414.8571
$ws.Range("H141").Value = 3096.647
$ws.Range("I141").Value = 3177.6875
$ws.Range("J141").Value = 1800
$ws.Range("K141").Value = 9533.0625
$ws.Range("L141").Value = 5400
$ws.Range("M141").Value = -4353.0625
$ws.Range("N141").Value = -15760
$ws.Range("M18").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5032.875
$ws.Range("I26").Value = 3894.7144
$ws.Range("K26").Value = 3894.7144
$ws.Range("M26").Value = -3564.7144
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 7000
$ws.Range("J29").Value = 13000
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 13000
$ws.Range("M29").Value = -6692
$ws.Range("N29").Value = -13616
$ws.Range("H32").Value = 2948.1592
$ws.Range("I32").Value = 2731.4048
$ws.Range("K32").Value = 2731.4048
$ws.Range("M32").Value = -2444.4048
$ws.Range("H74").Value = 1221.8334
$ws.Range("I74").Value = 1254.8
$ws.Range("K74").Value = 1254.8
$ws.Range("M74").Value = -380.8
$ws.Range("H77").Value = 1221.8334
$ws.Range("I77").Value = 1254.8
$ws.Range("K77").Value = 6274
$ws.Range("M77").Value = -1906
$ws.Range("H132").Value = 2255.7896
$ws.Range("I132").Value = 2255.7896
$ws.Range("K132").Value = 6767.3688
$ws.Range("M132").Value = -4237.3688

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 31443.75
$ws.Range("I82").Value = 18041.666
$ws.Range("K82").Value = 18041.666
$ws.Range("M82").Value = -17658.666
$ws.Range("H85").Value = 31443.75
$ws.Range("I85").Value = 18041.666
$ws.Range("K85").Value = 18041.666
$ws.Range("M85").Value = -16715.666
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H134").Value = 6086.25
$ws.Range("I134").Value = 6086.25
$ws.Range("K134").Value = 18258.75
$ws.Range("M134").Value = -15723.75
$ws.Range("N132").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("H51").Value = 22765
$ws.Range("J51").Value = 24900
$ws.Range("L51").Value = 24900
$ws.Range("N51").Value = -26372
$ws.Range("H61").Value = 22765
$ws.Range("J61").Value = 24900
$ws.Range("L61").Value = 24900
$ws.Range("N61").Value = -25596
$ws.Range("H62").Value = 3999
$ws.Range("J62").Value = 3999
$ws.Range("L62").Value = 3999
$ws.Range("N62").Value = -5247
$ws.Range("H65").Value = 3999
$ws.Range("J65").Value = 3999
$ws.Range("L65").Value = 19995
$ws.Range("N65").Value = -26235
$ws.Range("H99").Value = 3336.4443
$ws.Range("I99").Value = 3946.8572
$ws.Range("K99").Value = 3946.8572
$ws.Range("M99").Value = -2448.8572
$ws.Range("H107").Value = 889.5
$ws.Range("I107").Value = 815.5
$ws.Range("J107").Value = 945
$ws.Range("K107").Value = 815.5
$ws.Range("L107").Value = 945
$ws.Range("M107").Value = 1104.5
$ws.Range("N107").Value = -4785
$ws.Range("H126").Value = 3336.4443
$ws.Range("I126").Value = 3946.8572
$ws.Range("K126").Value = 11840.5716
$ws.Range("M126").Value = -9370.571599999999
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H132").Value = 2565
$ws.Range("I132").Value = 2565
$ws.Range("K132").Value = 7695
$ws.Range("M132").Value = -5165
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("N127").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99849.836
$ws.Range("J37").Value = 99849.836
$ws.Range("L37").Value = 299549.508
$ws.Range("N37").Value = -299773.508
$ws.Range("H47").Value = 926
$ws.Range("I47").Value = 675.3333
$ws.Range("J47").Value = 1302
$ws.Range("K47").Value = 2025.9999
$ws.Range("L47").Value = 3906
$ws.Range("M47").Value = -1594.9999
$ws.Range("N47").Value = -4768
$ws.Range("H50").Value = 981.75
$ws.Range("I50").Value = 416.5
$ws.Range("K50").Value = 1249.5
$ws.Range("M50").Value = -768.5
$ws.Range("H53").Value = 981.75
$ws.Range("I53").Value = 416.5
$ws.Range("K53").Value = 1249.5
$ws.Range("M53").Value = -768.5
$ws.Range("H55").Value = 2553.889
$ws.Range("J55").Value = 3019.2856
$ws.Range("L55").Value = 9057.856800000001
$ws.Range("N55").Value = -9411.856800000001
$ws.Range("H122").Value = 479.6
$ws.Range("J122").Value = 399.5
$ws.Range("L122").Value = 3595.5
$ws.Range("N122").Value = -8495.5
$ws.Range("H136").Value = 1999
$ws.Range("I136").Value = 1999
$ws.Range("K136").Value = 5997
$ws.Range("M136").Value = -897

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15000
$ws.Range("I46").Value = 8333.333000000001
$ws.Range("K46").Value = 8333.333000000001
$ws.Range("M46").Value = -8177.333000000001
$ws.Range("H102").Value = 2470.75
$ws.Range("I102").Value = 1635.5333
$ws.Range("K102").Value = 1635.5333
$ws.Range("M102").Value = -13.53330000000005
$ws.Range("H132").Value = 4183.8125
$ws.Range("I132").Value = 4129.4
$ws.Range("K132").Value = 12388.2
$ws.Range("M132").Value = -9858.199999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1538.6
$ws.Range("I16").Value = 1538.6
$ws.Range("K16").Value = 1538.6
$ws.Range("M16").Value = -1368.6
$ws.Range("H40").Value = 1685.2858
$ws.Range("I40").Value = 1655.2727
$ws.Range("J40").Value = 1795.3334
$ws.Range("K40").Value = 1655.2727
$ws.Range("L40").Value = 1795.3334
$ws.Range("M40").Value = -1519.2727
$ws.Range("N40").Value = -2067.3334
$ws.Range("H55").Value = 270.33334
$ws.Range("I55").Value = 277.72726
$ws.Range("J55").Value = 262.2
$ws.Range("K55").Value = 277.72726
$ws.Range("L55").Value = 262.2
$ws.Range("M55").Value = -104.72726
$ws.Range("N55").Value = -608.2
$ws.Range("H122").Value = 3402.2
$ws.Range("I122").Value = 3402.2
$ws.Range("K122").Value = 10206.6
$ws.Range("M122").Value = -7756.599999999999
$ws.Range("H132").Value = 3333.3333
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -17810
$ws.Range("H136").Value = 4072.4285
$ws.Range("I136").Value = 3504
$ws.Range("K136").Value = 10512
$ws.Range("M136").Value = -7962

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 34349.75
$ws.Range("J41").Value = 34349.75
$ws.Range("L41").Value = 34349.75
$ws.Range("N41").Value = -35129.75
$ws.Range("H122").Value = 2798.8125
$ws.Range("I122").Value = 2600.1538
$ws.Range("J122").Value = 3659.6667
$ws.Range("K122").Value = 7800.4614
$ws.Range("L122").Value = 10979.0001
$ws.Range("M122").Value = -5350.4614
$ws.Range("N122").Value = -15879.0001
$ws.Range("H126").Value = 984.2143
$ws.Range("J126").Value = 1329.8
$ws.Range("L126").Value = 3989.4
$ws.Range("N126").Value = -8929.4
$ws.Range("H132").Value = 1453.1875
$ws.Range("I132").Value = 1341
$ws.Range("K132").Value = 4023
$ws.Range("M132").Value = -1493
$ws.Range("H136").Value = 3370.027
$ws.Range("J136").Value = 2709.75
$ws.Range("L136").Value = 8129.25
$ws.Range("N136").Value = -13229.25
